$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report row was inserted at row 242 (week of 2021-10-20, "Zafiro
# rojo" / Primera), pushing the previously-existing rows 242-353 down to
# 243-354 (dimension grows from A1:R353 to A1:R354).
$ws.Range("A242").EntireRow.Insert()

$newRow = 242
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item($newRow, 3).Value = 'Maule'
$ws.Cells.Item($newRow, 4).Value = 44489
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100112002
$ws.Cells.Item($newRow, 7).Value = 'Pimiento'
$ws.Cells.Item($newRow, 8).Value = 'Zafiro rojo'
$ws.Cells.Item($newRow, 9).Value = 'Primera'
$ws.Cells.Item($newRow, 10).Value = 150
$ws.Cells.Item($newRow, 11).Value = 43000
$ws.Cells.Item($newRow, 12).Value = 43000
$ws.Cells.Item($newRow, 13).Value = 43000
$ws.Cells.Item($newRow, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item($newRow, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item($newRow, 16).Value = 2867
$ws.Cells.Item($newRow, 17).Value = 15
$ws.Cells.Item($newRow, 18).Value = 'Hortaliza'
